# Delete the "VW Polo 2022" row (row 36) from the "Data" sheet.
# This shifts all subsequent rows up by one, so what was row 37
# (Nissan Ariya 2022) becomes the new row 36, and so on, with the
# last row (previously 67 - Cupra Born 2022) disappearing and the
# used range shrinking from A1:J67 to A1:J66.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

$ws.Rows.Item(36).Delete()
